# Finished Week 13 logging
# Update Row 3 ("R") stats on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# OFF sheet
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 172
$wsOff.Range("C3").Value = 134
$wsOff.Range("D3").Value = 43
$wsOff.Range("E3").Value = 28

# DEF sheet
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 140
$wsDef.Range("C3").Value = 93
$wsDef.Range("D3").Value = 28
$wsDef.Range("E3").Value = 12
$wsDef.Range("F3").Value = 5
